$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header for column D and clear column E (the "Phone" column is removed entirely)
$ws.Range("D1").Value = "Folio No"

# Replace the SPV codes (AAA..EEE) with numeric folio numbers 1..5
$ws.Range("D2").Value = 1
$ws.Range("D3").Value = 2
$ws.Range("D4").Value = 3
$ws.Range("D5").Value = 4
$ws.Range("D6").Value = 5

# Remove column E entirely (Phone numbers column)
$ws.Range("E1:E6").Delete()

# Update selection to reflect the new state (no longer a full-column selection)
$ws.Range("D7").Select()
